# Apply cryptocurrency price/volume updates from the Fri May 12 21:59:56 UTC 2023
# GitHub Actions scraper run. Column D (Price) and E (Volume 1h) are refreshed
# for most rows; rows 36/37 additionally swap the Frax / TrustWalletToken
# entries (coin name, link, price, volume) to reflect their new rank order.
#
# Note: several "Price" values look numeric (e.g. "1.001") but must stay as
# literal text to match the original inlineStr cell type/formatting (the sheet
# also has genuinely dotted-thousands values like "26.799.56" that Excel would
# never auto-parse). Assigning such strings via .Value directly causes Excel's
# COM layer to auto-coerce them into floating point numbers, which would
# silently change the stored cell type/precision. We avoid that by writing a
# leading apostrophe (forcing literal text, exactly like a user typing '1.001
# into the cell) and then resetting the cell style to "Normal" to strip the
# transient quote-prefix style flag Excel attaches, so the only thing that
# changes versus the original file is the text content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.799.56'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.808.70'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'310.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'0.4302"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("D8").Value = "'0.3675"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.14%  '
$ws.Range("D9").Value = "'0.07195"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").Value = "'0.8607"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = "'20.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("D12").Value = '1.944.47'
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").Value = "'6.586"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("D14").Value = "'5.335"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").Value = "'0.06883"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = "'80.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008819"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = "'15.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '26.828.52'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = "'5.196"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.50%  '
$ws.Range("D23").Value = "'11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = '2.152.20'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").Value = "'152.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = "'1.863"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.33%  '
$ws.Range("D27").Value = "'18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("D28").Value = "'5.202"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.39%  '
$ws.Range("D29").Value = "'1.898"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.80%  '
$ws.Range("D30").Value = "'115.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("D31").Value = "'0.08922"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = "'0.7538"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("D33").Value = "'1.166"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.28%  '
$ws.Range("D34").Value = "'4.407"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").Value = "'2.773"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.05%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = "'1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = "'1.124"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.11%  '
$ws.Range("D38").Value = "'0.05190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").Value = "'0.01918"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").Value = "'0.5074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("D41").Value = "'0.1645"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("D42").Value = "'2.652"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = "'6.505"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.97%  '
$ws.Range("D44").Value = "'8.279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").Value = "'106.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = "'10.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = "'1.651"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.95%  '
$ws.Range("D49").Value = "'0.4564"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.06264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Value = "'1.793"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.29%  '
